{"js": "// 1) \"Optical flow\" paragraph: \"...still water.\" -> \"...still water with no reflection.\"\nconst body = context.document.body;\nconst flowResults = body.search(\"a clear sky, still water.\", { matchCase: true });\nflowResults.load(\"text\");\nawait context.sync();\n\nif (flowResults.items.length > 0) {\n  flowResults.items[0].insertText(\n    \"a clear sky, still water with no reflection.\",\n    \"Replace\"\n  );\n  await context.sync();\n}\n\n// 2) Insert a new paragraph about RANSAC vs. Hough Transform right after the\n//    question \"What are the advantages of RANSAC when compared with Hough\n//    Transform?\" (between the two blank BodyText paragraphs that follow it).\nconst questionResults = body.search(\n  \"What are the advantages of RANSAC when compared with Hough Transform?\",\n  { matchCase: true }\n);\nquestionResults.load(\"paragraphs\");\nawait context.sync();\n\nconst questionParagraph = questionResults.items[0].paragraphs.getFirst();\nquestionParagraph.load(\"text\");\nawait context.sync();\n\nconst firstBlankParagraph = questionParagraph.getNext();\nfirstBlankParagraph.load(\"text\");\nawait context.sync();\n\nfirstBlankParagraph.insertParagraph(\n  \"RANSAC detects outliers and excludes them from the calculations involved in predicting shapes/lines, while Hough Transform includes every datapoint.\",\n  \"After\"\n);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) \"Optical flow\" paragraph: \"...still water.\" -> \"...still water with no reflection.\"\n$d.Content.Find.Execute(\n    \"a clear sky, still water.\",\n    $false, $false, $false, $false, $false, $true, 1, $false,\n    \"a clear sky, still water with no reflection.\",\n    2\n)\n\n# 2) Insert a new paragraph about RANSAC vs. Hough Transform right after the\n#    question \"What are the advantages of RANSAC when compared with Hough\n#    Transform?\" (between the two blank BodyText paragraphs that follow it).\n$questionParagraph = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*advantages of RANSAC when compared with Hough Transform?*\") {\n        $questionParagraph = $p\n        break\n    }\n}\n\n$firstBlankParagraph = $questionParagraph.Next()\n$firstBlankParagraph.Range.InsertParagraphAfter()\n$newParagraph = $firstBlankParagraph.Next()\n$newParagraph.Range.Text = \"RANSAC detects outliers and excludes them from the calculations involved in predicting shapes/lines, while Hough Transform includes every datapoint.\"\n"}
